$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates — force text storage so values like "88.20" or
# "0.998" keep their exact original formatting instead of being
# auto-converted into numbers by Excel's type inference.
$dUpdates = @(
    @{ Addr = "D2"; Value = '43.524.59' },
    @{ Addr = "D3"; Value = '2.604.22' },
    @{ Addr = "D4"; Value = '0.998' },
    @{ Addr = "D5"; Value = '301.48' },
    @{ Addr = "D6"; Value = '96.35' },
    @{ Addr = "D7"; Value = '0.578' },
    @{ Addr = "D9"; Value = '0.558' },
    @{ Addr = "D10"; Value = '36.77' },
    @{ Addr = "D11"; Value = '0.0816' },
    @{ Addr = "D12"; Value = '7.82' },
    @{ Addr = "D13"; Value = '3.000.95' },
    @{ Addr = "D15"; Value = '2.603.79' },
    @{ Addr = "D16"; Value = '0.892' },
    @{ Addr = "D17"; Value = '14.36' },
    @{ Addr = "D18"; Value = '43.428.10' },
    @{ Addr = "D19"; Value = '6.68' },
    @{ Addr = "D20"; Value = '0.0₃0979' },
    @{ Addr = "D21"; Value = '12.33' },
    @{ Addr = "D22"; Value = '72.79' },
    @{ Addr = "D23"; Value = '265.91' },
    @{ Addr = "D25"; Value = '2.20' },
    @{ Addr = "D26"; Value = '29.28' },
    @{ Addr = "D28"; Value = '10.28' },
    @{ Addr = "D30"; Value = '37.66' },
    @{ Addr = "D31"; Value = '6.06' },
    @{ Addr = "D32"; Value = '3.61' },
    @{ Addr = "D33"; Value = '2.23' },
    @{ Addr = "D34"; Value = '151.57' },
    @{ Addr = "D36"; Value = '0.0813' },
    @{ Addr = "D37"; Value = '0.118' },
    @{ Addr = "D38"; Value = '24.41' },
    @{ Addr = "D40"; Value = '16.72' },
    @{ Addr = "D41"; Value = '3.58' },
    @{ Addr = "D42"; Value = '0.0314' },
    @{ Addr = "D43"; Value = '3.86' },
    @{ Addr = "D44"; Value = '2.043.28' },
    @{ Addr = "D45"; Value = '0.995' },
    @{ Addr = "D46"; Value = '88.20' },
    @{ Addr = "D47"; Value = '9.13' },
    @{ Addr = "D49"; Value = '2.853.19' },
    @{ Addr = "D50"; Value = '105.94' },
    @{ Addr = "D51"; Value = '0.191' }
)

foreach ($u in $dUpdates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

# Column E (Volume/1h % change) updates — values already contain
# non-numeric characters (%, spaces) so Excel stores them as text as-is.
$eUpdates = @(
    @{ Addr = "E2"; Value = '  -6.23%  ' },
    @{ Addr = "E3"; Value = '  +0.08%  ' },
    @{ Addr = "E4"; Value = '  -0.25%  ' },
    @{ Addr = "E5"; Value = '  -1.98%  ' },
    @{ Addr = "E6"; Value = '  -3.90%  ' },
    @{ Addr = "E7"; Value = '  -4.27%  ' },
    @{ Addr = "E8"; Value = '  -0.04%  ' },
    @{ Addr = "E9"; Value = '  -3.45%  ' },
    @{ Addr = "E10"; Value = '  -6.64%  ' },
    @{ Addr = "E11"; Value = '  -3.53%  ' },
    @{ Addr = "E12"; Value = '  -4.24%  ' },
    @{ Addr = "E13"; Value = '  -0.12%  ' },
    @{ Addr = "E14"; Value = '  +1.01%  ' },
    @{ Addr = "E15"; Value = '  -0.42%  ' },
    @{ Addr = "E16"; Value = '  -3.26%  ' },
    @{ Addr = "E17"; Value = '  -4.09%  ' },
    @{ Addr = "E18"; Value = '  -6.76%  ' },
    @{ Addr = "E19"; Value = '  -1.33%  ' },
    @{ Addr = "E20"; Value = '  -3.99%  ' },
    @{ Addr = "E21"; Value = '  -5.11%  ' },
    @{ Addr = "E22"; Value = '  +1.44%  ' },
    @{ Addr = "E23"; Value = '  -4.95%  ' },
    @{ Addr = "E24"; Value = '  -3.71%  ' },
    @{ Addr = "E25"; Value = '  +1.22%  ' },
    @{ Addr = "E26"; Value = '  +1.84%  ' },
    @{ Addr = "E27"; Value = '  +0.32%  ' },
    @{ Addr = "E28"; Value = '  -3.34%  ' },
    @{ Addr = "E29"; Value = '  -1.23%  ' },
    @{ Addr = "E30"; Value = '  -3.90%  ' },
    @{ Addr = "E31"; Value = '  -5.09%  ' },
    @{ Addr = "E32"; Value = '  -1.40%  ' },
    @{ Addr = "E33"; Value = '  +0.67%  ' },
    @{ Addr = "E34"; Value = '  -0.08%  ' },
    @{ Addr = "E35"; Value = '  -2.23%  ' },
    @{ Addr = "E36"; Value = '  -3.39%  ' },
    @{ Addr = "E37"; Value = '  -4.31%  ' },
    @{ Addr = "E38"; Value = '  +5.80%  ' },
    @{ Addr = "E39"; Value = '  -1.43%  ' },
    @{ Addr = "E40"; Value = '  +3.49%  ' },
    @{ Addr = "E41"; Value = '  -1.43%  ' },
    @{ Addr = "E42"; Value = '  -5.51%  ' },
    @{ Addr = "E43"; Value = '  -5.41%  ' },
    @{ Addr = "E44"; Value = '  -4.12%  ' },
    @{ Addr = "E45"; Value = '  -0.33%  ' },
    @{ Addr = "E46"; Value = '  -5.29%  ' },
    @{ Addr = "E47"; Value = '  -4.49%  ' },
    @{ Addr = "E48"; Value = '  +3.60%  ' },
    @{ Addr = "E49"; Value = '  -0.25%  ' },
    @{ Addr = "E50"; Value = '  -3.29%  ' },
    @{ Addr = "E51"; Value = '  -4.75%  ' }
)

foreach ($u in $eUpdates) {
    $ws.Range($u.Addr).Value = $u.Value
}

Write-Host "Updated $($dUpdates.Count) price cells and $($eUpdates.Count) volume cells"
